# Add a new worksheet "2025-03-23" at the end of the workbook containing the
# 2025-03-23 price-summary data (same shape as the existing daily sheets).
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2025-03-23"

# Header row
$ws.Range("A1").Value = "Match"
$ws.Range("B1").Value = "Seat Type"
$ws.Range("C1").Value = "Min_Price"
$ws.Range("D1").Value = "Avg_Price"
$ws.Range("E1").Value = "Ticket_Count"

# Match the bold / centered header look used on the other daily sheets.
$ws.Range("A1:E1").Font.Bold = $true
$ws.Range("A1:E1").HorizontalAlignment = -4108
$ws.Range("A1:E1").VerticalAlignment = -4160
$ws.Range("A1:E1").Borders.LineStyle = 1

$ws.Range("A2").Value = 'Arsenal - Fulham'
$ws.Range("B2").Value = 'Shortside Upper'
$ws.Range("C2").Value = 84
$ws.Range("D2").Value = 202
$ws.Range("E2").Value = 58
$ws.Range("A3").Value = 'Arsenal - Fulham'
$ws.Range("B3").Value = 'Longside Lower'
$ws.Range("C3").Value = 98
$ws.Range("D3").Value = 295
$ws.Range("E3").Value = 111
$ws.Range("A4").Value = 'Arsenal - Fulham'
$ws.Range("B4").Value = 'Shortside Lower'
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = 235
$ws.Range("E4").Value = 66
$ws.Range("A5").Value = 'Arsenal - Fulham'
$ws.Range("B5").Value = 'Longside Upper'
$ws.Range("C5").Value = 110
$ws.Range("D5").Value = 260
$ws.Range("E5").Value = 54
$ws.Range("A6").Value = 'Arsenal - Fulham'
$ws.Range("B6").Value = 'Club Level'
$ws.Range("C6").Value = 213
$ws.Range("D6").Value = 611
$ws.Range("E6").Value = 16
$ws.Range("A7").Value = 'Arsenal - Fulham'
$ws.Range("B7").Value = 'VIP & Executive Box'
$ws.Range("C7").Value = 995
$ws.Range("D7").Value = 995
$ws.Range("E7").Value = 1
$ws.Range("A8").Value = 'Arsenal - Brentford'
$ws.Range("B8").Value = 'Longside Upper'
$ws.Range("C8").Value = 83
$ws.Range("D8").Value = 256
$ws.Range("E8").Value = 57
$ws.Range("A9").Value = 'Arsenal - Brentford'
$ws.Range("B9").Value = 'Shortside Lower'
$ws.Range("C9").Value = 95
$ws.Range("D9").Value = 257
$ws.Range("E9").Value = 60
$ws.Range("A10").Value = 'Arsenal - Brentford'
$ws.Range("B10").Value = 'Shortside Upper'
$ws.Range("C10").Value = 105
$ws.Range("D10").Value = 246
$ws.Range("E10").Value = 50
$ws.Range("A11").Value = 'Arsenal - Brentford'
$ws.Range("B11").Value = 'Longside Lower'
$ws.Range("C11").Value = 107
$ws.Range("D11").Value = 332
$ws.Range("E11").Value = 78
$ws.Range("A12").Value = 'Arsenal - Brentford'
$ws.Range("B12").Value = 'Away Section'
$ws.Range("C12").Value = 138
$ws.Range("D12").Value = 138
$ws.Range("E12").Value = 1
$ws.Range("A13").Value = 'Arsenal - Brentford'
$ws.Range("B13").Value = 'Club Level'
$ws.Range("C13").Value = 233
$ws.Range("D13").Value = 855
$ws.Range("E13").Value = 17
$ws.Range("A14").Value = 'Arsenal - Crystal Palace'
$ws.Range("B14").Value = 'Shortside Upper'
$ws.Range("C14").Value = 119
$ws.Range("D14").Value = 237
$ws.Range("E14").Value = 86
$ws.Range("A15").Value = 'Arsenal - Crystal Palace'
$ws.Range("B15").Value = 'Longside Lower'
$ws.Range("C15").Value = 128
$ws.Range("D15").Value = 345
$ws.Range("E15").Value = 108
$ws.Range("A16").Value = 'Arsenal - Crystal Palace'
$ws.Range("B16").Value = 'Shortside Lower'
$ws.Range("C16").Value = 133
$ws.Range("D16").Value = 274
$ws.Range("E16").Value = 84
$ws.Range("A17").Value = 'Arsenal - Crystal Palace'
$ws.Range("B17").Value = 'Longside Upper'
$ws.Range("C17").Value = 155
$ws.Range("D17").Value = 269
$ws.Range("E17").Value = 82
$ws.Range("A18").Value = 'Arsenal - Crystal Palace'
$ws.Range("B18").Value = 'Away Section'
$ws.Range("C18").Value = 261
$ws.Range("D18").Value = 261
$ws.Range("E18").Value = 1
$ws.Range("A19").Value = 'Arsenal - Crystal Palace'
$ws.Range("B19").Value = 'Club Level'
$ws.Range("C19").Value = 299
$ws.Range("D19").Value = 716
$ws.Range("E19").Value = 18
$ws.Range("A20").Value = 'Arsenal v Real Madrid : Champions League 2024-2025'
$ws.Range("B20").Value = 'Shortside Upper'
$ws.Range("C20").Value = 592
$ws.Range("D20").Value = 750
$ws.Range("E20").Value = 38
$ws.Range("A21").Value = 'Arsenal v Real Madrid : Champions League 2024-2025'
$ws.Range("B21").Value = 'Longside Upper'
$ws.Range("C21").Value = 603
$ws.Range("D21").Value = 871
$ws.Range("E21").Value = 59
$ws.Range("A22").Value = 'Arsenal v Real Madrid : Champions League 2024-2025'
$ws.Range("B22").Value = 'Shortside Lower'
$ws.Range("C22").Value = 662
$ws.Range("D22").Value = 850
$ws.Range("E22").Value = 36
$ws.Range("A23").Value = 'Arsenal v Real Madrid : Champions League 2024-2025'
$ws.Range("B23").Value = 'Longside Lower'
$ws.Range("C23").Value = 663
$ws.Range("D23").Value = 1201
$ws.Range("E23").Value = 92
$ws.Range("A24").Value = 'Arsenal v Real Madrid : Champions League 2024-2025'
$ws.Range("B24").Value = 'Club Level'
$ws.Range("C24").Value = 884
$ws.Range("D24").Value = 1638
$ws.Range("E24").Value = 29
$ws.Range("A25").Value = 'Arsenal v Real Madrid : Champions League 2024-2025'
$ws.Range("B25").Value = 'Away Section'
$ws.Range("C25").Value = 1658
$ws.Range("D25").Value = 2488
$ws.Range("E25").Value = 4
$ws.Range("A26").Value = 'Arsenal v Real Madrid : Champions League 2024-2025'
$ws.Range("B26").Value = 'VIP & Executive Box'
$ws.Range("C26").Value = 3870
$ws.Range("D26").Value = 4699
$ws.Range("E26").Value = 2
